# "Improve introspect texts fixes #30"
#
# 1. Refresh the "datetimeFigureOut" date placeholders (slide master +
#    every slide layout) from 10/19/2020 to 11/12/2020.
# 2. On slide 1, the "TextBox 20" shape drops its first line
#    ("Get Authorization Token [ITI-71]") - it auto-shrinks afterwards
#    because the box uses <a:spAutoFit/>.

$p = $ppt.ActivePresentation
$cr = [char]13

function Update-DatePlaceholder($shapes) {
    $cnt = $shapes.Count
    for ($i = 1; $i -le $cnt; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq "10/19/2020") {
                    $tr.Text = "11/12/2020"
                }
            }
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lyt = $layouts.Item($li)
    Update-DatePlaceholder $lyt.Shapes
}

# Slide 1: trim the first line out of the "TextBox 20" callout.
$s = $p.Slides.Item(1)
$shapes = $s.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.Name -eq "TextBox 20") {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf($cr)
        if ($idx -ge 0) {
            $tr.Text = $full.Substring($idx + 1)
        }
    }
}
